$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 210.3077
$ws.Range("I12").Value = 210.3077
$ws.Range("K12").Value = 210.3077
$ws.Range("M12").Value = -40.30770000000001
$ws.Range("H33").Value = 153.46153
$ws.Range("I33").Value = 161.5
$ws.Range("K33").Value = 161.5
$ws.Range("M33").Value = 67.5
$ws.Range("H98").Value = 4340.6
$ws.Range("I98").Value = 3990.3125
$ws.Range("J98").Value = 5741.75
$ws.Range("K98").Value = 3990.3125
$ws.Range("L98").Value = 5741.75
$ws.Range("M98").Value = -2492.3125
$ws.Range("N98").Value = -8737.75
$ws.Range("H113").Value = 4038.75
$ws.Range("I113").Value = 3159.6
$ws.Range("J113").Value = 5504
$ws.Range("K113").Value = 3159.6
$ws.Range("L113").Value = 5504
$ws.Range("M113").Value = 94.40000000000009
$ws.Range("N113").Value = -12012
$ws.Range("H122").Value = 4340.6
$ws.Range("I122").Value = 3990.3125
$ws.Range("J122").Value = 5741.75
$ws.Range("K122").Value = 11970.9375
$ws.Range("L122").Value = 17225.25
$ws.Range("M122").Value = -9520.9375
$ws.Range("N122").Value = -22125.25
$ws.Range("H132").Value = 6442.44
$ws.Range("I132").Value = 6850.478
$ws.Range("K132").Value = 20551.434
$ws.Range("M132").Value = -18021.434
$ws.Range("H136").Value = 85236
$ws.Range("J136").Value = 84999.5
$ws.Range("L136").Value = 84999.5
$ws.Range("N136").Value = -95199.5
$ws.Range("H137").Value = 1045980.25
$ws.Range("I137").Value = 1251670.2
$ws.Range("J137").Value = 17530.5
$ws.Range("K137").Value = 3755010.6
$ws.Range("L137").Value = 52591.5
$ws.Range("M137").Value = -3752460.6
$ws.Range("N137").Value = -57691.5
$ws.Range("H138").Value = 2370.4644
$ws.Range("I138").Value = 898.2727
$ws.Range("J138").Value = 3323.0588
$ws.Range("K138").Value = 2694.8181
$ws.Range("L138").Value = 9969.1764
$ws.Range("M138").Value = 2445.1819
$ws.Range("N138").Value = -20249.1764

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 144333.97
$ws.Range("I74").Value = 207241.78
$ws.Range("K74").Value = 207241.78
$ws.Range("M74").Value = -206367.78
$ws.Range("H77").Value = 144333.97
$ws.Range("I77").Value = 207241.78
$ws.Range("K77").Value = 1036208.9
$ws.Range("M77").Value = -1031840.9
$ws.Range("H122").Value = 3897.8333
$ws.Range("I122").Value = 3977.4
$ws.Range("K122").Value = 11932.2
$ws.Range("M122").Value = -9482.200000000001
$ws.Range("H132").Value = 6799.125
$ws.Range("I132").Value = 8499.5
$ws.Range("J132").Value = 6232.3335
$ws.Range("K132").Value = 25498.5
$ws.Range("L132").Value = 18697.0005
$ws.Range("M132").Value = -22968.5
$ws.Range("N132").Value = -23757.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 598
$ws.Range("J80").Value = 846
$ws.Range("L80").Value = 846
$ws.Range("N80").Value = -2842
$ws.Range("H83").Value = 598
$ws.Range("J83").Value = 846
$ws.Range("L83").Value = 4230
$ws.Range("N83").Value = -14214
$ws.Range("H134").Value = 3321.8572
$ws.Range("J134").Value = 4564.25
$ws.Range("L134").Value = 13692.75
$ws.Range("N134").Value = -18762.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2376.1667
$ws.Range("I16").Value = 2064.25
$ws.Range("K16").Value = 2064.25
$ws.Range("M16").Value = -1777.25
$ws.Range("H31").Value = 3576818
$ws.Range("I31").Value = 3862.6924
$ws.Range("K31").Value = 3862.6924
$ws.Range("M31").Value = -3567.6924
$ws.Range("H34").Value = 3576818
$ws.Range("I34").Value = 3862.6924
$ws.Range("K34").Value = 3862.6924
$ws.Range("M34").Value = -3660.6924
$ws.Range("H50").Value = 60417.8
$ws.Range("J50").Value = 60417.8
$ws.Range("L50").Value = 60417.8
$ws.Range("N50").Value = -61667.8
$ws.Range("H86").Value = 5312
$ws.Range("I86").Value = 3899.75
$ws.Range("J86").Value = 6253.5
$ws.Range("K86").Value = 3899.75
$ws.Range("L86").Value = 6253.5
$ws.Range("M86").Value = -2776.75
$ws.Range("N86").Value = -8499.5
$ws.Range("H89").Value = 5312
$ws.Range("I89").Value = 3899.75
$ws.Range("J89").Value = 6253.5
$ws.Range("K89").Value = 19498.75
$ws.Range("L89").Value = 31267.5
$ws.Range("M89").Value = -13882.75
$ws.Range("N89").Value = -42499.5
$ws.Range("H113").Value = 2376.1667
$ws.Range("I113").Value = 2064.25
$ws.Range("K113").Value = 2064.25
$ws.Range("M113").Value = 105.75
$ws.Range("H122").Value = 507.375
$ws.Range("I122").Value = 491.25
$ws.Range("K122").Value = 1473.75
$ws.Range("M122").Value = 976.25
$ws.Range("H132").Value = 3528.5625
$ws.Range("I132").Value = 3391.3
$ws.Range("K132").Value = 10173.9
$ws.Range("M132").Value = -7643.900000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 895.1429000000001
$ws.Range("I9").Value = 158.66667
$ws.Range("J9").Value = 1447.5
$ws.Range("K9").Value = 476.00001
$ws.Range("L9").Value = 4342.5
$ws.Range("M9").Value = -252.00001
$ws.Range("N9").Value = -4790.5
$ws.Range("H74").Value = 26873.625
$ws.Range("J74").Value = 27496
$ws.Range("L74").Value = 82488
$ws.Range("N74").Value = -84610
$ws.Range("H77").Value = 26873.625
$ws.Range("J77").Value = 27496
$ws.Range("L77").Value = 247464
$ws.Range("N77").Value = -258072
$ws.Range("H137").Value = 1957.6364
$ws.Range("I137").Value = 1142.3334
$ws.Range("J137").Value = 2936
$ws.Range("K137").Value = 3427.0002
$ws.Range("L137").Value = 8808
$ws.Range("M137").Value = 1672.9998
$ws.Range("N137").Value = -19008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6944.1387
$ws.Range("I102").Value = 998.6667
$ws.Range("K102").Value = 998.6667
$ws.Range("M102").Value = 623.3333
$ws.Range("H107").Value = 10244.2
$ws.Range("I107").Value = 2589.6
$ws.Range("K107").Value = 2589.6
$ws.Range("M107").Value = -669.5999999999999
$ws.Range("H113").Value = 2436.125
$ws.Range("I113").Value = 2453
$ws.Range("K113").Value = 2453
$ws.Range("M113").Value = -283
$ws.Range("H126").Value = 9035
$ws.Range("I126").Value = 2360
$ws.Range("J126").Value = 13206.875
$ws.Range("K126").Value = 7080
$ws.Range("L126").Value = 39620.625
$ws.Range("M126").Value = -4610
$ws.Range("N126").Value = -44560.625
$ws.Range("H132").Value = 2569.9312
$ws.Range("I132").Value = 2334.6
$ws.Range("J132").Value = 2693.7896
$ws.Range("K132").Value = 7003.799999999999
$ws.Range("L132").Value = 8081.3688
$ws.Range("M132").Value = -4473.799999999999
$ws.Range("N132").Value = -13141.3688
$ws.Range("J141").Value = 24000
$ws.Range("L141").Value = 24000
$ws.Range("N141").Value = -34360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 19762.154
$ws.Range("I40").Value = 21159
$ws.Range("K40").Value = 21159
$ws.Range("M40").Value = -21023
$ws.Range("H122").Value = 8888.190000000001
$ws.Range("J122").Value = 9833.272000000001
$ws.Range("L122").Value = 29499.816
$ws.Range("N122").Value = -34399.81600000001
$ws.Range("H132").Value = 5841.2593
$ws.Range("I132").Value = 3790.1333
$ws.Range("J132").Value = 8405.166999999999
$ws.Range("K132").Value = 11370.3999
$ws.Range("L132").Value = 25215.501
$ws.Range("M132").Value = -8840.3999
$ws.Range("N132").Value = -30275.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7666.1113
$ws.Range("I81").Value = 10499.5
$ws.Range("K81").Value = 20999
$ws.Range("M81").Value = -19938
$ws.Range("H84").Value = 7666.1113
$ws.Range("I84").Value = 10499.5
$ws.Range("K84").Value = 104995
$ws.Range("M84").Value = -99691
$ws.Range("H122").Value = 13158549
$ws.Range("J122").Value = 125000500
$ws.Range("L122").Value = 375001500
$ws.Range("N122").Value = -375006400
$ws.Range("H132").Value = 1623.4445
$ws.Range("I132").Value = 1549.6666
$ws.Range("K132").Value = 4648.9998
$ws.Range("M132").Value = -2118.9998
$ws.Range("H140").Value = 57886.668
$ws.Range("J140").Value = 57886.668
$ws.Range("L140").Value = 57886.668
$ws.Range("N140").Value = -68246.66800000001
